# "Add files via upload" — data corrections to Gegevens_avIncidentenApp.xlsx
#
# 1) Sheet "Oplossingen": remove the row for OplossingID=2 (Excel row 3);
#    all subsequent rows shift up one, which EntireRow delete does natively.
# 2) Sheet "Handelingen": row for IncidentID was missing its ID (A4 blank);
#    fill it in with 3.
# Plus matching selection changes left behind by the editor on each sheet.

$wb = $excel.ActiveWorkbook

$wsOplossingen = $wb.Worksheets.Item("Oplossingen")
$wsOplossingen.Rows(3).Delete()
$wsOplossingen.Range("A3:XFD3").Select()

$wsHandelingen = $wb.Worksheets.Item("Handelingen")
$wsHandelingen.Range("A4").Value = 3
$wsHandelingen.Range("C13").Select()
